$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) The shared string "upb401" (cell A2) should also match uppercase searches,
#    so normalize it to "uPB401".
$ws.Range("A2").Value = "uPB401"

# 2) Row height of the data rows (2-12) grew slightly (18.75 -> 19.5pt).
$ws.Range("A2:B12").RowHeight = 19.5

# 3) Re-colour the data rows' font to an explicit black (was inheriting the
#    theme colour) and make sure the header row keeps/gets the same
#    (border+font) formatting that was already used by the data rows.
$ws.Range("A2:B12").Font.Color = -16777216
$ws.Range("A1:B1").Font.Name = "Calibri"
